$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: re-style to match the "end of block" look (border row, like row 4) ---
$ws.Range("A4:E4").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)

# --- Add rows 8 and 9: a new two-row block, styled like rows 2/3 (start of block + continuation) ---
$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

$ws.Range("B3:E3").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)

$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 21.6

# --- Fill in the new data (column-major so shared strings land in source order) ---
$ws.Range("A8").Value = "SCRIPT/T01P02A/um1105.ssb"

$ws.Range("B8").Value = 179
$ws.Range("B9").Value = 182

$ws.Range("C8").Value = ' Our sense of smell\''s acute.\nWe\''re real sensitive about it.'
$ws.Range("C9").Value = ' But the stink is mostly gone.\nThat\''s a good thing.'

$ws.Range("D8").Value = ' У нас великолепный нюх.\nМы очень чувствительны к запахам.'
$ws.Range("D9").Value = ' Но вонь почти исчезла. Это\nпрекрасно.'

$ws.Range("E8").Value = ' Ô îàò âåìéëïìåðîúê îýö.\nÍú ïœåîû œôâòóâéóåìûîú ë èàðàöàí.'
$ws.Range("E9").Value = ' Îï âïîû ðïœóé éòœåèìà. Üóï\nðñåëñàòîï.'

# --- Selection matches the committed state ---
[void]$ws.Range("D5").Select()
